$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H11").Value = 67.57143000000001
$ws_ALC.Range("I11").Value = 67.57143000000001
$ws_ALC.Range("K11").Value = 67.57143000000001
$ws_ALC.Range("M11").Value = 72.42856999999999
$ws_ALC.Range("H18").Value = 687.375
$ws_ALC.Range("I18").Value = 366.66666
$ws_ALC.Range("J18").Value = 1649.5
$ws_ALC.Range("K18").Value = 366.66666
$ws_ALC.Range("L18").Value = 1649.5
$ws_ALC.Range("M18").Value = -82.66665999999998
$ws_ALC.Range("N18").Value = -2217.5
$ws_ALC.Range("H40").Value = 10818.454
$ws_ALC.Range("J40").Value = 11125
$ws_ALC.Range("L40").Value = 11125
$ws_ALC.Range("N40").Value = -11475
$ws_ALC.Range("H74").Value = 11894.357
$ws_ALC.Range("I74").Value = 10886.23
$ws_ALC.Range("K74").Value = 10886.23
$ws_ALC.Range("M74").Value = -9950.23
$ws_ALC.Range("H77").Value = 11894.357
$ws_ALC.Range("I77").Value = 10886.23
$ws_ALC.Range("K77").Value = 54431.14999999999
$ws_ALC.Range("M77").Value = -49751.14999999999
$ws_ALC.Range("H132").Value = 4372.533
$ws_ALC.Range("I132").Value = 4328.0303
$ws_ALC.Range("K132").Value = 12984.0909
$ws_ALC.Range("M132").Value = -10454.0909
$ws_ALC.Range("H133").Value = 58334.332
$ws_ALC.Range("J133").Value = 58334.332
$ws_ALC.Range("L133").Value = 58334.332
$ws_ALC.Range("N133").Value = -68454.33199999999
$ws_ALC.Range("H137").Value = 4983.115
$ws_ALC.Range("I137").Value = 5665.8335
$ws_ALC.Range("J137").Value = 4397.9287
$ws_ALC.Range("K137").Value = 16997.5005
$ws_ALC.Range("L137").Value = 13193.7861
$ws_ALC.Range("M137").Value = -14447.5005
$ws_ALC.Range("N137").Value = -18293.7861
$ws_ALC.Range("H138").Value = 4180.5303
$ws_ALC.Range("I138").Value = 1564.5333
$ws_ALC.Range("J138").Value = 4949.9414
$ws_ALC.Range("K138").Value = 4693.5999
$ws_ALC.Range("L138").Value = 14849.8242
$ws_ALC.Range("M138").Value = 446.4000999999998
$ws_ALC.Range("N138").Value = -25129.8242
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H61").Value = 4946.769
$ws_ARM.Range("I61").Value = 2330
$ws_ARM.Range("K61").Value = 2330
$ws_ARM.Range("M61").Value = -2118
$ws_ARM.Range("H74").Value = 28129.455
$ws_ARM.Range("I74").Value = 68520.25
$ws_ARM.Range("J74").Value = 5049
$ws_ARM.Range("K74").Value = 68520.25
$ws_ARM.Range("L74").Value = 5049
$ws_ARM.Range("M74").Value = -67646.25
$ws_ARM.Range("N74").Value = -6797
$ws_ARM.Range("H77").Value = 28129.455
$ws_ARM.Range("I77").Value = 68520.25
$ws_ARM.Range("J77").Value = 5049
$ws_ARM.Range("K77").Value = 342601.25
$ws_ARM.Range("L77").Value = 25245
$ws_ARM.Range("M77").Value = -338233.25
$ws_ARM.Range("N77").Value = -33981
$ws_ARM.Range("H110").Value = 162588.19
$ws_ARM.Range("I110").Value = 193688.92
$ws_ARM.Range("J110").Value = 864.4
$ws_ARM.Range("K110").Value = 193688.92
$ws_ARM.Range("L110").Value = 864.4
$ws_ARM.Range("M110").Value = -191643.92
$ws_ARM.Range("N110").Value = -4954.4
$ws_ARM.Range("H132").Value = 5020.936
$ws_ARM.Range("I132").Value = 2610.8333
$ws_ARM.Range("K132").Value = 7832.499899999999
$ws_ARM.Range("M132").Value = -5302.499899999999
$ws_ARM.Range("H136").Value = 4946.769
$ws_ARM.Range("I136").Value = 2330
$ws_ARM.Range("K136").Value = 6990
$ws_ARM.Range("M136").Value = -4440
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H94").Value = 1041.6875
$ws_BSM.Range("I94").Value = 748.9167
$ws_BSM.Range("K94").Value = 748.9167
$ws_BSM.Range("M94").Value = -297.9167
$ws_BSM.Range("H132").Value = 62500
$ws_BSM.Range("J132").Value = 62500
$ws_BSM.Range("L132").Value = 62500
$ws_BSM.Range("N132").Value = -72620
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H22").Value = 493.125
$ws_CRP.Range("I22").Value = 492.14285
$ws_CRP.Range("K22").Value = 492.14285
$ws_CRP.Range("M22").Value = -142.14285
$ws_CRP.Range("H86").Value = 39999.668
$ws_CRP.Range("I86").Value = 9999
$ws_CRP.Range("J86").Value = 55000
$ws_CRP.Range("K86").Value = 9999
$ws_CRP.Range("L86").Value = 55000
$ws_CRP.Range("M86").Value = -8876
$ws_CRP.Range("N86").Value = -57246
$ws_CRP.Range("H89").Value = 39999.668
$ws_CRP.Range("I89").Value = 9999
$ws_CRP.Range("J89").Value = 55000
$ws_CRP.Range("K89").Value = 49995
$ws_CRP.Range("L89").Value = 275000
$ws_CRP.Range("M89").Value = -44379
$ws_CRP.Range("N89").Value = -286232
$ws_CRP.Range("H132").Value = 3304.775
$ws_CRP.Range("I132").Value = 2399.9
$ws_CRP.Range("K132").Value = 7199.700000000001
$ws_CRP.Range("M132").Value = -4669.700000000001
$ws_CRP.Range("H134").Value = 3829.1707
$ws_CRP.Range("I134").Value = 3109.375
$ws_CRP.Range("J134").Value = 6388.4443
$ws_CRP.Range("K134").Value = 9328.125
$ws_CRP.Range("L134").Value = 19165.3329
$ws_CRP.Range("M134").Value = -6793.125
$ws_CRP.Range("N134").Value = -24235.3329
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H6").Value = 3499.6667
$ws_GSM.Range("J6").Value = 3499.6667
$ws_GSM.Range("L6").Value = 3499.6667
$ws_GSM.Range("N6").Value = -3725.6667
$ws_GSM.Range("H16").Value = 3499.6667
$ws_GSM.Range("J16").Value = 3499.6667
$ws_GSM.Range("L16").Value = 3499.6667
$ws_GSM.Range("N16").Value = -3999.6667
$ws_GSM.Range("H20").Value = 28666.666
$ws_GSM.Range("J20").Value = 28666.666
$ws_GSM.Range("L20").Value = 28666.666
$ws_GSM.Range("N20").Value = -29156.666
$ws_GSM.Range("H24").Value = 3353333.2
$ws_GSM.Range("I24").Value = 60000
$ws_GSM.Range("K24").Value = 60000
$ws_GSM.Range("M24").Value = -59827
$ws_GSM.Range("H113").Value = 404360.72
$ws_GSM.Range("J113").Value = 13997.5
$ws_GSM.Range("L113").Value = 13997.5
$ws_GSM.Range("N113").Value = -18337.5
$ws_GSM.Range("H125").Value = 75000
$ws_GSM.Range("J125").Value = 75000
$ws_GSM.Range("L125").Value = 75000
$ws_GSM.Range("N125").Value = -79920
$ws_GSM.Range("H126").Value = 200003540
$ws_GSM.Range("I126").Value = 333336000
$ws_GSM.Range("J126").Value = 4849.5
$ws_GSM.Range("K126").Value = 1000008000
$ws_GSM.Range("L126").Value = 14548.5
$ws_GSM.Range("M126").Value = -1000005530
$ws_GSM.Range("N126").Value = -19488.5
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H7").Value = 674612.9399999999
$ws_LTW.Range("I7").Value = 916799.6
$ws_LTW.Range("K7").Value = 916799.6
$ws_LTW.Range("M7").Value = -916687.6
$ws_LTW.Range("H16").Value = 3227.4614
$ws_LTW.Range("I16").Value = 2945.2222
$ws_LTW.Range("K16").Value = 2945.2222
$ws_LTW.Range("M16").Value = -2775.2222
$ws_LTW.Range("H40").Value = 835920.9399999999
$ws_LTW.Range("I40").Value = 1001655.1
$ws_LTW.Range("K40").Value = 1001655.1
$ws_LTW.Range("M40").Value = -1001519.1
$ws_LTW.Range("H46").Value = 5470.5293
$ws_LTW.Range("I46").Value = 4149.9
$ws_LTW.Range("K46").Value = 4149.9
$ws_LTW.Range("M46").Value = -3961.9
$ws_LTW.Range("H61").Value = 4337.5757
$ws_LTW.Range("I61").Value = 3248.5
$ws_LTW.Range("J61").Value = 6013.077
$ws_LTW.Range("K61").Value = 3248.5
$ws_LTW.Range("L61").Value = 6013.077
$ws_LTW.Range("M61").Value = -3046.5
$ws_LTW.Range("N61").Value = -6417.077
$ws_LTW.Range("H100").Value = 251651.25
$ws_LTW.Range("I100").Value = 251651.25
$ws_LTW.Range("J100").Value = 0
$ws_LTW.Range("K100").Value = 251651.25
$ws_LTW.Range("L100").Value = 0
$ws_LTW.Range("M100").Value = -251110.25
$ws_LTW.Range("N100").ClearContents()
$ws_LTW.Range("H113").Value = 4337.5757
$ws_LTW.Range("I113").Value = 3248.5
$ws_LTW.Range("J113").Value = 6013.077
$ws_LTW.Range("K113").Value = 3248.5
$ws_LTW.Range("L113").Value = 6013.077
$ws_LTW.Range("M113").Value = -1078.5
$ws_LTW.Range("N113").Value = -10353.077
$ws_LTW.Range("H126").Value = 674612.9399999999
$ws_LTW.Range("I126").Value = 916799.6
$ws_LTW.Range("K126").Value = 2750398.8
$ws_LTW.Range("M126").Value = -2747928.8
$ws_LTW.Range("H132").Value = 6822.5454
$ws_LTW.Range("J132").Value = 8256.125
$ws_LTW.Range("L132").Value = 24768.375
$ws_LTW.Range("N132").Value = -29828.375
$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H6").Value = 5000
$ws_WVR.Range("J6").Value = 5000
$ws_WVR.Range("L6").Value = 5000
$ws_WVR.Range("N6").Value = -5230
$ws_WVR.Range("H46").Value = 0
$ws_WVR.Range("J46").Value = 0
$ws_WVR.Range("L46").Value = 0
$ws_WVR.Range("N46").ClearContents()
$ws_WVR.Range("H49").Value = 10000
$ws_WVR.Range("J49").Value = 10000
$ws_WVR.Range("L49").Value = 10000
$ws_WVR.Range("N49").Value = -10460
$ws_WVR.Range("H100").Value = 297.73685
$ws_WVR.Range("I100").Value = 271.2857
$ws_WVR.Range("K100").Value = 542.5714
$ws_WVR.Range("M100").Value = -1.57140000000004
$ws_WVR.Range("H107").Value = 31203.47
$ws_WVR.Range("I107").Value = 40504.23
$ws_WVR.Range("J107").Value = 976
$ws_WVR.Range("K107").Value = 121512.69
$ws_WVR.Range("L107").Value = 2928
$ws_WVR.Range("M107").Value = -119592.69
$ws_WVR.Range("N107").Value = -6768
$ws_WVR.Range("H132").Value = 4899.067
$ws_WVR.Range("I132").Value = 4623.8335
$ws_WVR.Range("K132").Value = 13871.5005
$ws_WVR.Range("M132").Value = -11341.5005
$ws_WVR.Range("H134").Value = 0
$ws_WVR.Range("J134").Value = 0
$ws_WVR.Range("L134").Value = 0
$ws_WVR.Range("N134").ClearContents()
